$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours worked for Ian Stolz (row 6) and Bernard Swanepoel (row 10)
$ws.Range("C6").Value = 3.5
$ws.Range("C10").Value = 3.5

# Update the last selection to C21
$ws.Range("C21").Select()
